# Fix exploration time formula (#13)
# Update Avg_Agent_Step_Time (G), Avg_Experiment_Time (H),
# Std_Agent_Step_Time (M), and Std_Experiment_Time (N) columns
# with the corrected simulation results.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 4.918226689999999
$ws.Range("H2").Value = 1232.07125056
$ws.Range("M2").Value = 1.777316944179572
$ws.Range("N2").Value = 963.7916867193834

$ws.Range("G3").Value = 2.49977837
$ws.Range("H3").Value = 835.5034485
$ws.Range("M3").Value = 0.790676197487009
$ws.Range("N3").Value = 484.5582803624869

$ws.Range("G4").Value = 2.57720064
$ws.Range("H4").Value = 331.60441235
$ws.Range("M4").Value = 0.785403263238125
$ws.Range("N4").Value = 216.1385229431625

$ws.Range("G5").Value = 1.12823029
$ws.Range("H5").Value = 192.20037335
$ws.Range("M5").Value = 0.3516925354369894
$ws.Range("N5").Value = 102.7755886657826

$ws.Range("G6").Value = 1.29440734
$ws.Range("H6").Value = 86.81829558999999
$ws.Range("M6").Value = 0.3599765913815329
$ws.Range("N6").Value = 50.43125049219724

$ws.Range("G7").Value = 0.46392292
$ws.Range("H7").Value = 41.55343851
$ws.Range("M7").Value = 0.1758878481641225
$ws.Range("N7").Value = 25.48424586805145

$ws.Range("G8").Value = 0.8823541
$ws.Range("H8").Value = 41.7716212
$ws.Range("M8").Value = 0.2374836300065988
$ws.Range("N8").Value = 23.09148883162552

$ws.Range("G9").Value = 0.24769386
$ws.Range("H9").Value = 15.31183321
$ws.Range("M9").Value = 0.104901358710742
$ws.Range("N9").Value = 9.616772294894949

$ws.Range("G10").Value = 0.65727415
$ws.Range("H10").Value = 24.30724414
$ws.Range("M10").Value = 0.1601921226960032
$ws.Range("N10").Value = 11.6996009233499

$ws.Range("G11").Value = 0.14931174
$ws.Range("H11").Value = 7.236461
$ws.Range("M11").Value = 0.07534440205013702
$ws.Range("N11").Value = 5.104913695057046

$ws.Range("G12").Value = 0.52116752
$ws.Range("H12").Value = 16.20027588
$ws.Range("M12").Value = 0.1384449157439339
$ws.Range("N12").Value = 7.750529262724177

$ws.Range("G13").Value = 0.09075055999999998
$ws.Range("H13").Value = 3.52285935
$ws.Range("M13").Value = 0.05167540004156577
$ws.Range("N13").Value = 2.821168194457838

